$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A6').Value = '<small class="nav-text text-muted me-auto" data-bs-toggle="tooltip" data-bs-placement="bottom" title="">0.2.4</small>'

$ws.Range('A15').Value = '<ul class="dropdown-menu" aria-labelledby="dropdown-articles"><li><a class="dropdown-item" href="../articles/data_mapping.html">Mapping variables to outputs</a></li>'
$ws.Range('A16').Value = '<li><a class="dropdown-item" href="../articles/table_templates.html">CRVS tables</a></li>'
$ws.Range('A17').Value = '</ul></li>'
$ws.Range('A18').Value = '</ul><ul class="navbar-nav"><li class="nav-item"><form class="form-inline" role="search">'
$ws.Range('A19').Value = '<input class="form-control" type="search" name="search-input" id="search-input" autocomplete="off" aria-label="Search site" placeholder="Search for" data-search-index="../search.json"></form></li>'
$ws.Range('A20').Value = '<li class="nav-item dropdown">'
$ws.Range('A21').Value = '<button class="nav-link dropdown-toggle" type="button" id="dropdown-lightswitch" data-bs-toggle="dropdown" aria-expanded="false" aria-haspopup="true" aria-label="Light switch"><span class="fa fa-sun"></span></button>'
$ws.Range('A22').Value = '<ul class="dropdown-menu dropdown-menu-end" aria-labelledby="dropdown-lightswitch"><li><button class="dropdown-item" data-bs-theme-value="light"><span class="fa fa-sun"></span> Light</button></li>'
$ws.Range('A23').Value = '<li><button class="dropdown-item" data-bs-theme-value="dark"><span class="fa fa-moon"></span> Dark</button></li>'
$ws.Range('A24').Value = '<li><button class="dropdown-item" data-bs-theme-value="auto"><span class="fa fa-adjust"></span> Auto</button></li>'
$ws.Range('A25').Value = '</ul></li>'
$ws.Range('A26').Value = '</ul></div>'
$ws.Range('A27').Value = '</div>'
$ws.Range('A28').Value = '</nav><div class="container template-reference-topic">'
$ws.Range('A29').Value = '<div class="row">'
$ws.Range('A30').Value = '<main id="main" class="col-md-9"><div class="page-header">'
$ws.Range('A31').Value = '<h1>Convert CSV Files to an Excel Workbook</h1>'
$ws.Range('A32').Value = '<div class="d-none name"><code>convert_csv_xlsx.Rd</code></div>'
$ws.Range('A33').Value = '</div>'
$ws.Range('A34').Value = '<div class="ref-description section level2">'
$ws.Range('A35').Value = '<p>This function reads all CSV files in a specified directory and writes their contents to separate sheets in a single Excel workbook.</p>'
$ws.Range('A36').Value = '</div>'
$ws.Range('A37').Value = '<div class="section level2">'
$ws.Range('A38').Value = '<h2 id="ref-usage">Usage<a class="anchor" aria-label="anchor" href="#ref-usage"></a></h2>'
$ws.Range('A39').Value = '<div class="sourceCode"><pre class="sourceCode r"><code><span><span class="fu">convert_csv_xlsx</span><span class="op">(</span>input_path <span class="op">=</span> <span class="st">"."</span>, output_path <span class="op">=</span> <span class="st">"output.xlsx"</span><span class="op">)</span></span></code></pre></div>'
$ws.Range('A40').Value = '</div>'
$ws.Range('A41').Value = '<div class="section level2">'
$ws.Range('A42').Value = '<h2 id="arguments">Arguments<a class="anchor" aria-label="anchor" href="#arguments"></a></h2>'
$ws.Range('A43').Value = '<dl><dt id="arg-output-path">output_path<a class="anchor" aria-label="anchor" href="#arg-output-path"></a></dt>'
$ws.Range('A44').Value = '<dd><p>A character string specifying the directory and file name to write the xlsx to.</p></dd>'
$ws.Range('A45').Value = '<dt id="arg-path">path<a class="anchor" aria-label="anchor" href="#arg-path"></a></dt>'
$ws.Range('A46').Value = '<dd><p>A character string specifying the directory containing the CSV files. Defaults to the current working directory (".").</p></dd>'
$ws.Range('A47').Value = '</dl></div>'
$ws.Range('A48').Value = '<div class="section level2">'
$ws.Range('A49').Value = '<h2 id="value">Value<a class="anchor" aria-label="anchor" href="#value"></a></h2>'
$ws.Range('A50').Value = '<p>An Excel file named "output.xlsx" containing the contents of the CSV files.</p>'
$ws.Range('A51').Value = '</div>'
$ws.Range('A52').Value = '<div class="section level2">'
$ws.Range('A53').Value = '<h2 id="ref-examples">Examples<a class="anchor" aria-label="anchor" href="#ref-examples"></a></h2>'
$ws.Range('A54').Value = '<div class="sourceCode"><pre class="sourceCode r"><code><span class="r-in"><span><span class="co"># Convert CSV files in the current directory to an Excel workbook</span></span></span>'
$ws.Range('A55').Value = '<span class="r-in"><span><span class="fu"><a href="csv_to_excel.html">csv_to_excel</a></span><span class="op">(</span>path <span class="op">=</span> <span class="st">"."</span><span class="op">)</span></span></span>'
$ws.Range('A56').Value = '<span class="r-wrn co"><span class="r-pr">#&gt;</span> <span class="warning">Warning: </span>There are no .csv files in this directory</span>'
$ws.Range('A57').Value = '<span class="r-out co"><span class="r-pr">#&gt;</span> NULL</span>'
$ws.Range('A58').Value = '<span class="r-in"><span></span></span>'
$ws.Range('A59').Value = '<span class="r-in"><span><span class="co"># Convert CSV files in a specified directory to an Excel workbook</span></span></span>'
$ws.Range('A60').Value = '<span class="r-in"><span><span class="fu"><a href="csv_to_excel.html">csv_to_excel</a></span><span class="op">(</span>input_path <span class="op">=</span> <span class="st">"path/to/directory"</span>, output_path <span class="op">=</span> <span class="st">"path/with/file/output.xlsx"</span><span class="op">)</span></span></span>'
$ws.Range('A61').Value = '<span class="r-err co"><span class="r-pr">#&gt;</span> <span class="error">Error in csv_to_excel(input_path = "path/to/directory", output_path = "path/with/file/output.xlsx"):</span> unused arguments (input_path = "path/to/directory", output_path = "path/with/file/output.xlsx")</span>'
$ws.Range('A62').Value = '</code></pre></div>'
$ws.Range('A63').Value = '</div>'
$ws.Range('A64').Value = '</main><aside class="col-md-3"><nav id="toc" aria-label="Table of contents"><h2>On this page</h2>'
$ws.Range('A65').Value = '</nav></aside></div>'
$ws.Range('A66').Value = '<footer><div class="pkgdown-footer-left">'
$ws.Range('A67').Value = '<p>Developed by Tesfaye Belay, Pamela Kakande, Rachel Shipsey, Liam Beardsmore.</p>'
$ws.Range('A68').Value = '</div>'
$ws.Range('A69').Value = '<div class="pkgdown-footer-right">'
$ws.Range('A70').Value = '<p>Site built with <a href="https://pkgdown.r-lib.org/" class="external-link">pkgdown</a> 2.1.0.</p>'
$ws.Range('A71').Value = '</div>'
$ws.Range('A72').Value = '</footer></div>'
$ws.Range('A73').Value = '</body></html>'

$ws.Range("A74:A76").EntireRow.Delete() | Out-Null